# Scheduled runner update: refresh currentAveragePrice-derived profit figures
# across the Seraph_Profits crafting-leve sheets (market price pull + recalculated
# NQ/HQ totals and profit deltas).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2000
$ws.Range("J17").Value = 2000
$ws.Range("L17").Value = 6000
$ws.Range("N17").Value = -6336
# Row 18
$ws.Range("H18").Value = 2116.5557
$ws.Range("I18").Value = 2116.5557
$ws.Range("K18").Value = 2116.5557
$ws.Range("M18").Value = -1832.5557
# Row 40
$ws.Range("H40").Value = 2264.4285
$ws.Range("I40").Value = 1975
$ws.Range("K40").Value = 1975
$ws.Range("M40").Value = -1800
# Row 55
$ws.Range("H55").Value = 241.125
$ws.Range("J55").Value = 325.8
$ws.Range("L55").Value = 325.8
$ws.Range("N55").Value = -753.8
# Row 93
$ws.Range("H93").Value = 20000
$ws.Range("J93").Value = 20000
$ws.Range("L93").Value = 20000
$ws.Range("N93").Value = -24992
# Row 132
$ws.Range("H132").Value = 1242.7727
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = $null
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 1744.2363
$ws.Range("I74").Value = 1168.8125
$ws.Range("K74").Value = 1168.8125
$ws.Range("M74").Value = -294.8125
# Row 77
$ws.Range("H77").Value = 1744.2363
$ws.Range("I77").Value = 1168.8125
$ws.Range("K77").Value = 5844.0625
$ws.Range("M77").Value = -1476.0625
# Row 110
$ws.Range("H110").Value = 6408
$ws.Range("I110").Value = 7948
$ws.Range("J110").Value = 248
$ws.Range("K110").Value = 7948
$ws.Range("L110").Value = 248
$ws.Range("M110").Value = -5903
$ws.Range("N110").Value = -4338
# Row 122
$ws.Range("H122").Value = 2318.3076
$ws.Range("I122").Value = 1682.1904
$ws.Range("K122").Value = 5046.5712
$ws.Range("M122").Value = -2596.5712
# Row 132
$ws.Range("H132").Value = 2504.5715
$ws.Range("I132").Value = 2504.5715
$ws.Range("K132").Value = 7513.7145
$ws.Range("M132").Value = -4983.7145
$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 624.6
$ws.Range("I64").Value = 498.375
$ws.Range("J64").Value = 768.8570999999999
$ws.Range("K64").Value = 498.375
$ws.Range("L64").Value = 768.8570999999999
$ws.Range("M64").Value = -273.375
$ws.Range("N64").Value = -1218.8571
# Row 67
$ws.Range("H67").Value = 624.6
$ws.Range("I67").Value = 498.375
$ws.Range("J67").Value = 768.8570999999999
$ws.Range("K67").Value = 498.375
$ws.Range("L67").Value = 768.8570999999999
$ws.Range("M67").Value = 281.625
$ws.Range("N67").Value = -2328.8571
# Row 86
$ws.Range("H86").Value = 1017.0769
$ws.Range("I86").Value = 1086.3334
$ws.Range("J86").Value = 861.25
$ws.Range("K86").Value = 1086.3334
$ws.Range("L86").Value = 861.25
$ws.Range("M86").Value = 36.66660000000002
$ws.Range("N86").Value = -3107.25
# Row 89
$ws.Range("H89").Value = 1017.0769
$ws.Range("I89").Value = 1086.3334
$ws.Range("J89").Value = 861.25
$ws.Range("K89").Value = 5431.666999999999
$ws.Range("L89").Value = 4306.25
$ws.Range("M89").Value = 184.3330000000005
$ws.Range("N89").Value = -15538.25
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 289.64285
$ws.Range("I22").Value = 223.55556
$ws.Range("K22").Value = 223.55556
$ws.Range("M22").Value = 126.44444
# Row 122
$ws.Range("H122").Value = 3450.56
$ws.Range("I122").Value = 3218.6365
$ws.Range("K122").Value = 9655.9095
$ws.Range("M122").Value = -7205.9095
# Row 132
$ws.Range("H132").Value = 6983.3335
$ws.Range("J132").Value = 6991
$ws.Range("L132").Value = 20973
$ws.Range("N132").Value = -26033
$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 1290
$ws.Range("J22").Value = 300
$ws.Range("L22").Value = 900
$ws.Range("N22").Value = -1238
# Row 27
$ws.Range("H27").Value = 1290
$ws.Range("J27").Value = 300
$ws.Range("L27").Value = 900
$ws.Range("N27").Value = -1104
# Row 32
$ws.Range("H32").Value = 6980772
$ws.Range("J32").Value = 11633645
$ws.Range("L32").Value = 34900935
$ws.Range("N32").Value = -34901501
# Row 114
$ws.Range("H114").Value = 260
$ws.Range("I114").Value = 223.66667
$ws.Range("J114").Value = 287.25
$ws.Range("K114").Value = 671.00001
$ws.Range("L114").Value = 861.75
$ws.Range("M114").Value = 2582.99999
$ws.Range("N114").Value = -7369.75
# Row 118
$ws.Range("H118").Value = 900
$ws.Range("I118").Value = 900
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 2700
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -1457
$ws.Range("N118").Value = $null
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2465.6875
$ws.Range("I80").Value = 1989.8422
$ws.Range("J80").Value = 3161.1538
$ws.Range("K80").Value = 1989.8422
$ws.Range("L80").Value = 3161.1538
$ws.Range("M80").Value = -991.8422
$ws.Range("N80").Value = -5157.1538
# Row 83
$ws.Range("H83").Value = 2465.6875
$ws.Range("I83").Value = 1989.8422
$ws.Range("J83").Value = 3161.1538
$ws.Range("K83").Value = 9949.210999999999
$ws.Range("L83").Value = 15805.769
$ws.Range("M83").Value = -4957.210999999999
$ws.Range("N83").Value = -25789.769
# Row 102
$ws.Range("H102").Value = 1719.619
$ws.Range("I102").Value = 407.53333
$ws.Range("K102").Value = 407.53333
$ws.Range("M102").Value = 1214.46667
# Row 122
$ws.Range("H122").Value = 64979.688
$ws.Range("I122").Value = 2291.1428
$ws.Range("K122").Value = 6873.428400000001
$ws.Range("M122").Value = -4423.428400000001
# Row 132
$ws.Range("H132").Value = 1362.8889
$ws.Range("I132").Value = 1191.8636
$ws.Range("K132").Value = 3575.5908
$ws.Range("M132").Value = -1045.5908
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 1595.3334
$ws.Range("I40").Value = 1650.75
$ws.Range("K40").Value = 1650.75
$ws.Range("M40").Value = -1514.75
# Row 55
$ws.Range("H55").Value = 410.85715
$ws.Range("J55").Value = 489
$ws.Range("L55").Value = 489
$ws.Range("N55").Value = -835
# Row 68
$ws.Range("H68").Value = 2120.5
$ws.Range("I68").Value = 1867.3334
$ws.Range("K68").Value = 1867.3334
$ws.Range("M68").Value = -1118.3334
# Row 71
$ws.Range("H71").Value = 2120.5
$ws.Range("I71").Value = 1867.3334
$ws.Range("K71").Value = 9336.666999999999
$ws.Range("M71").Value = -5592.666999999999
# Row 122
$ws.Range("H122").Value = 4833
$ws.Range("I122").Value = 4833
$ws.Range("K122").Value = 14499
$ws.Range("M122").Value = -12049
# Row 131
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").Value = $null
# Row 132
$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470
# Row 136
$ws.Range("H136").Value = 2399
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = $null
$ws = $wb.Worksheets.Item("WVR")
# Row 47
$ws.Range("H47").Value = 31229.75
$ws.Range("J47").Value = 31229.75
$ws.Range("L47").Value = 31229.75
$ws.Range("N47").Value = -32373.75
# Row 113
$ws.Range("H113").Value = 319.3846
$ws.Range("I113").Value = 367.33334
$ws.Range("K113").Value = 1102.00002
$ws.Range("M113").Value = 1067.99998
# Row 122
$ws.Range("H122").Value = 3672.2144
$ws.Range("I122").Value = 4150.9165
$ws.Range("K122").Value = 12452.7495
$ws.Range("M122").Value = -10002.7495
# Row 126
$ws.Range("H126").Value = 2678.4
$ws.Range("I126").Value = 2244.2727
$ws.Range("K126").Value = 6732.8181
$ws.Range("M126").Value = -4262.8181
# Row 132
$ws.Range("H132").Value = 1570.4
$ws.Range("J132").Value = 3000
$ws.Range("L132").Value = 9000
$ws.Range("N132").Value = -14060
# Row 136
$ws.Range("H136").Value = 8995
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 8995
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 26985
$ws.Range("M136").Value = $null
$ws.Range("N136").Value = -32085
